# Presentacion.pptx — "Agregado de menu (sin funcionalidad). Actualizacion de presentacion."
#
# This script reproduces, via PowerPoint COM automation:
#   1. Slide 1 ("menu" shapes being nudged horizontally) — five shapes get a
#      new horizontal (x / Left) offset while y/width/height stay put.
#   2. Slide 4 — a trailing run split ("Planes segun Necesidades
#      Especificas:" / " " / "Empresarial o Familiar") gets collapsed so the
#      last two runs become a single run " Empresarial o Familiar".
#
# Note on point values below: Shape.Left/Top are expressed in points, while
# the OOXML stores EMU (1 pt = 12700 EMU). The COM layer quantizes the point
# value through a float32 round-trip before converting to EMU, so a few of
# the literals here carry extra fractional digits chosen so that, after that
# quantization, they land exactly on the target EMU value from the diff
# (plain target_emu/12700 can land one EMU short because of the float32
# rounding). Each target is annotated with the EMU value it reproduces.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 1 — shift the five "menu" shapes horizontally.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Shape 1: "1 Titulo" — x 611560 -> 687760 EMU
$s1.Shapes.Item(1).Left = 54.15435028076172

# Shape 2: "2 Subtitulo" — x 899592 -> 902668 EMU
$s1.Shapes.Item(2).Left = 71.07625579833984

# Shape 3: "Picture 3" — x 2555776 -> 2833989 EMU
$s1.Shapes.Item(3).Left = 223.14877319335938

# Shape 4: "2 Subtitulo" (bottom caption) — x 827584 -> 866664 EMU
$s1.Shapes.Item(4).Left = 68.24129486083984

# Shape 5: "1 Titulo" (bottom "Panic Dial Button") — x 763960 -> 687760 EMU
$s1.Shapes.Item(5).Left = 54.15435028076172

# ---------------------------------------------------------------------
# 2) Slide 4 — merge the trailing " " + "Empresarial o Familiar" runs.
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# Locate the bullet textbox by content rather than a hard-coded index.
$sh = $null
for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $cand = $s4.Shapes.Item($i)
    if ($cand.HasTextFrame -and $cand.TextFrame.TextRange.Text -like "*Empresarial o Familiar*") {
        $sh = $cand
        break
    }
}

$tr = $sh.TextFrame.TextRange

# Find the paragraph containing the target phrase (paragraph 4 in the
# original deck: "Planes segun Necesidades Especificas: Empresarial o
# Familiar").
$para = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $cand = $tr.Paragraphs($i)
    if ($cand.Text -like "*Empresarial o Familiar*") {
        $para = $cand
        break
    }
}

# Grab the exact sub-range covering the last two underlying runs
# (" " and "Empresarial o Familiar") and rewrite it as one run's text;
# reassigning .Text across both source runs collapses them into a single
# run in the saved XML.
$mergeStart = $para.Text.IndexOf(" Empresarial o Familiar") + 1
$mergeLen = $para.Length - $mergeStart + 1
$sub = $para.Characters($mergeStart, $mergeLen)
$sub.Text = " Empresarial o Familiar"
